# Applies the "Updated symbol list" data refresh to the cryptos sheet.
# Every changed cell keeps its original inline/shared-string (text) type,
# even though many values look numeric (prices, percentages, hour digit),
# so NumberFormat is forced to "@" (Text) before each assignment to stop
# Excel from silently re-typing the cell as a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '291.11'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '-3.43%'
$ws.Range('G2').NumberFormat = '@'
$ws.Range('G2').Value = '8'
# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '30.66'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '-6.58%'
$ws.Range('G3').NumberFormat = '@'
$ws.Range('G3').Value = '8'
# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '4.955'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '0.20%'
$ws.Range('G4').NumberFormat = '@'
$ws.Range('G4').Value = '8'
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.07232'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '-6.80%'
$ws.Range('G5').NumberFormat = '@'
$ws.Range('G5').Value = '8'
# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.790'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '-8.47%'
$ws.Range('G6').NumberFormat = '@'
$ws.Range('G6').Value = '8'
# Row 7
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '-2.04%'
$ws.Range('G7').NumberFormat = '@'
$ws.Range('G7').Value = '8'
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.761'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '-1.06%'
$ws.Range('G8').NumberFormat = '@'
$ws.Range('G8').Value = '8'
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.8969'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '-2.54%'
$ws.Range('G9').NumberFormat = '@'
$ws.Range('G9').Value = '8'
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.1657'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '-6.56%'
$ws.Range('G10').NumberFormat = '@'
$ws.Range('G10').Value = '8'
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07711'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '-2.40%'
$ws.Range('G11').NumberFormat = '@'
$ws.Range('G11').Value = '8'
# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08079'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-6.27%'
$ws.Range('G12').NumberFormat = '@'
$ws.Range('G12').Value = '8'
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.03038'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-3.85%'
$ws.Range('G13').NumberFormat = '@'
$ws.Range('G13').Value = '8'
# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.1003'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '0.10%'
$ws.Range('G14').NumberFormat = '@'
$ws.Range('G14').Value = '8'
# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.001506'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '-0.89%'
$ws.Range('G15').NumberFormat = '@'
$ws.Range('G15').Value = '8'
# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.005863'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '1.28%'
$ws.Range('G16').NumberFormat = '@'
$ws.Range('G16').Value = '8'
# Row 17
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'UpBots'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.007492'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-0.12%'
$ws.Range('G17').NumberFormat = '@'
$ws.Range('G17').Value = '8'
# Row 18
$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'LEO'
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.474'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '0.34%'
$ws.Range('G18').NumberFormat = '@'
$ws.Range('G18').Value = '8'
# Row 19
$ws.Range('B19').NumberFormat = '@'
$ws.Range('B19').Value = 'BTSEToken'
$ws.Range('C19').NumberFormat = '@'
$ws.Range('C19').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.084'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '-3.26%'
$ws.Range('G19').NumberFormat = '@'
$ws.Range('G19').Value = '8'
# Row 20
$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = 'BitpandaEcosystemToken'
$ws.Range('C20').NumberFormat = '@'
$ws.Range('C20').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.3316'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-0.75%'
$ws.Range('G20').NumberFormat = '@'
$ws.Range('G20').Value = '8'
# Row 21
$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'ProBitToken'
$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.1289'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '-2.35%'
$ws.Range('G21').NumberFormat = '@'
$ws.Range('G21').Value = '8'
# Row 22
$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'MCDex'
$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.033'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '-6.11%'
$ws.Range('G22').NumberFormat = '@'
$ws.Range('G22').Value = '8'
# Row 23
$ws.Range('B23').NumberFormat = '@'
$ws.Range('B23').Value = 'ZBToken'
$ws.Range('C23').NumberFormat = '@'
$ws.Range('C23').Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.2251'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '13.07%'
$ws.Range('G23').NumberFormat = '@'
$ws.Range('G23').Value = '8'
# Row 24
$ws.Range('B24').NumberFormat = '@'
$ws.Range('B24').Value = 'CoinExToken'
$ws.Range('C24').NumberFormat = '@'
$ws.Range('C24').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.04510'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-0.81%'
$ws.Range('G24').NumberFormat = '@'
$ws.Range('G24').Value = '8'
# Row 25
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'BitKan'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.001214'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '-1.02%'
$ws.Range('G25').NumberFormat = '@'
$ws.Range('G25').Value = '8'
# Row 26
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'HotbitToken'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.004014'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '-9.33%'
$ws.Range('G26').NumberFormat = '@'
$ws.Range('G26').Value = '8'
# Row 27
$ws.Range('B27').NumberFormat = '@'
$ws.Range('B27').Value = 'NitroEx'
$ws.Range('C27').NumberFormat = '@'
$ws.Range('C27').Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0001250'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '-0.07%'
$ws.Range('G27').NumberFormat = '@'
$ws.Range('G27').Value = '8'
# Row 28
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'Spectre.aiUtilityToken'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut'
$ws.Range('G28').NumberFormat = '@'
$ws.Range('G28').Value = '8'
# Row 29
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'LegolasExchange'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo'
$ws.Range('G29').NumberFormat = '@'
$ws.Range('G29').Value = '8'
# Row 30
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'BitZToken'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz'
$ws.Range('G30').NumberFormat = '@'
$ws.Range('G30').Value = '8'
# Row 31
$ws.Range('B31').NumberFormat = '@'
$ws.Range('B31').Value = 'Birake'
$ws.Range('C31').NumberFormat = '@'
$ws.Range('C31').Value = 'https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir'
$ws.Range('G31').NumberFormat = '@'
$ws.Range('G31').Value = '8'
# Row 32
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'NashExchange'
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex'
$ws.Range('G32').NumberFormat = '@'
$ws.Range('G32').Value = '8'
# Row 33
$ws.Range('B33').NumberFormat = '@'
$ws.Range('B33').Value = 'AAXToken'
$ws.Range('C33').NumberFormat = '@'
$ws.Range('C33').Value = 'https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab'
$ws.Range('G33').NumberFormat = '@'
$ws.Range('G33').Value = '8'
# Row 34
$ws.Range('B34').NumberFormat = '@'
$ws.Range('B34').Value = 'CenX'
$ws.Range('C34').NumberFormat = '@'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V4XJUvLQb+cenx-cenx'
$ws.Range('G34').NumberFormat = '@'
$ws.Range('G34').Value = '8'
# Row 35
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'BNIXToken'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/n194X9uHp+bnixtoken-bnix'
$ws.Range('G35').NumberFormat = '@'
$ws.Range('G35').Value = '8'
# Row 36
$ws.Range('G36').NumberFormat = '@'
$ws.Range('G36').Value = '8'
# Row 37
$ws.Range('G37').NumberFormat = '@'
$ws.Range('G37').Value = '8'
# Row 38
$ws.Range('G38').NumberFormat = '@'
$ws.Range('G38').Value = '8'
# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01602'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '-6.04%'
$ws.Range('G39').NumberFormat = '@'
$ws.Range('G39').Value = '8'
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.04411'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '-6.47%'
$ws.Range('G40').NumberFormat = '@'
$ws.Range('G40').Value = '8'
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.007265'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '-3.56%'
$ws.Range('G41').NumberFormat = '@'
$ws.Range('G41').Value = '8'
# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1309'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '-3.29%'
$ws.Range('G42').NumberFormat = '@'
$ws.Range('G42').Value = '8'
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.007682'
$ws.Range('G43').NumberFormat = '@'
$ws.Range('G43').Value = '8'
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.001901'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '-18.86%'
$ws.Range('G44').NumberFormat = '@'
$ws.Range('G44').Value = '8'
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.009218'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '-12.13%'
$ws.Range('G45').NumberFormat = '@'
$ws.Range('G45').Value = '8'
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00005946'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '-5.23%'
$ws.Range('G46').NumberFormat = '@'
$ws.Range('G46').Value = '8'
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000750'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '-0.03%'
$ws.Range('G47').NumberFormat = '@'
$ws.Range('G47').Value = '8'
# Row 48
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '172.73%'
$ws.Range('G48').NumberFormat = '@'
$ws.Range('G48').Value = '8'
# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.003001'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-3.27%'
$ws.Range('G49').NumberFormat = '@'
$ws.Range('G49').Value = '8'
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002101'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '-0.03%'
$ws.Range('G50').NumberFormat = '@'
$ws.Range('G50').Value = '8'
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0002001'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '-0.03%'
$ws.Range('G51').NumberFormat = '@'
$ws.Range('G51').Value = '8'
